$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting existing rows 23:126 down to 24:127
$ws.Rows("23:23").Insert()

# Populate the newly inserted row 23 with the new record's data
$ws.Cells.Item(23, 1).Value = 4
$ws.Cells.Item(23, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(23, 3).Value = "Los Lagos"
$ws.Cells.Item(23, 4).Value = 44764
$ws.Cells.Item(23, 5).Value = 10
$ws.Cells.Item(23, 6).Value = 100112022
$ws.Cells.Item(23, 7).Value = "Arveja Verde"
$ws.Cells.Item(23, 8).Value = "Perfection"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 70
$ws.Cells.Item(23, 11).Value = 43000
$ws.Cells.Item(23, 12).Value = 43000
$ws.Cells.Item(23, 13).Value = 43000
$ws.Cells.Item(23, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(23, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(23, 16).Value = 1720
$ws.Cells.Item(23, 17).Value = 25
$ws.Cells.Item(23, 18).Value = "Hortaliza"
